# Actualización automática 2025-08-20 12:25:09
#
# A new advisor/client row ("ALTAMIRANO VILLAVICENCIO JUAN ALEJANDRO") is
# inserted at row 4 (alphabetically between "ALTAMIRANO ARIAS LUCIA
# ELIZABETH" and the old "ALVAREZ SAAVEDRA EDWIN GEOVANNY") on both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. Every row from the old
# row 4 onward shifts down by one, and the trailing summary row (counts
# "N de 56" on VENTAS POR GRUPO) becomes "N de 57" to reflect the new
# total of 57 clients.

$wb = $excel.ActiveWorkbook

$asesor = "LINDAO ZUÑIGA BRYAN JOSE"
$nuevoCliente = "ALTAMIRANO VILLAVICENCIO JUAN ALEJANDRO"

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" — columns A:R, data rows 2-58, summary 58->59
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row at 4, pushing the old row 4 ("ALVAREZ SAAVEDRA ...")
# and everything below it down by one.
$ws1.Rows.Item(4).Insert()

$ws1.Cells.Item(4, 1).Value2 = $asesor
$ws1.Cells.Item(4, 2).Value2 = $nuevoCliente
for ($col = 3; $col -le 18; $col++) {
    # Rows.Item(4).Insert() already pushed the style (s="2", currency
    # number format) down from the old row 4 into the new blank row 4,
    # so only the value needs to be written here.
    $ws1.Cells.Item(4, $col).Value2 = 0
}

# The trailing "N de 56" summary row is now row 59 — bump the count to 57.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(59, $col)
    $cell.Value2 = $cell.Value2 -replace "de 56", "de 57"
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" — columns A:G, data rows 2-58, summary 58->59
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(4).Insert()

$ws2.Cells.Item(4, 1).Value2 = $asesor
$ws2.Cells.Item(4, 2).Value2 = $nuevoCliente
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(4, $col).Value2 = 0
}

# Row 59's totals are plain numbers (not formulas) and are unaffected by
# the inserted all-zero row, so no further changes are needed there.
